# "div last news in frontpage" — append a new data row (row 37) to the
# "Remarque - Problematique" tracking sheet, mirroring the layout of the
# row above it (row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (borders, wrap, fonts, row height source) of the
# last existing data row onto the new row before filling in values.
$ws.Range("A36:H36").Copy()
$ws.Range("A37:H37").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(37).RowHeight = 82.5

# Fill in the new "problème" entry.
$ws.Range("B37").Value = "construire nouveau sujet a partir de article"                 # problème rencontré
$ws.Range("D37").Value = "routing , default text field value"                          # solution definitive
$ws.Range("E37").Value = 42110                                                          # date (2015-04-16)
$ws.Range("H37").Value = "http://stackoverflow.com/questions/13916001/set-default-value-of-symfony-2-form-field-in-twig"

# Turn the source link in column H into a real hyperlink, matching the
# bold+underlined "Lien hypertexte" look used by the other rows.
$ws.Hyperlinks.Add($ws.Cells.Item(37,8), "http://stackoverflow.com/questions/13916001/set-default-value-of-symfony-2-form-field-in-twig")
$ws.Range("H37").Font.Bold = $true

# Extend the visible selection to cover the newly added row.
$ws.Range("A5:H37").Select() | Out-Null

Write-Host "done"
